# Code v1.1 updated after running all the tests. Minor modifications made.
#
# Updates the "Currency" sheet's INR-converted values (B2:B9) with refreshed
# conversion rates, then restores the view state: "Currency" becomes the
# active/selected tab (with A2:B9 highlighted) while "Gold Price" loses its
# previous tab-selected state (with A2:D3 highlighted instead).

$wb = $excel.ActiveWorkbook

$wsCurrency = $wb.Worksheets.Item("Currency")
$wsGold     = $wb.Worksheets.Item("Gold Price")

# Refreshed currency conversion rates.
$wsCurrency.Range("B2").Value = "73.0675"
$wsCurrency.Range("B3").Value = "88.4809"
$wsCurrency.Range("B4").Value = "103.0862"
$wsCurrency.Range("B5").Value = "56.3058"
$wsCurrency.Range("B6").Value = "0.6664"
$wsCurrency.Range("B7").Value = "55.1203"
$wsCurrency.Range("B8").Value = "11.4227"
$wsCurrency.Range("B9").Value = "2.6434"

# Leave a selection on the "Gold Price" sheet before switching away from it,
# matching the saved selection the diff shows for that sheet.
$wsGold.Range("A2:D3").Select()

# "Currency" becomes the active tab/sheet.
$wsCurrency.Activate()
$wsCurrency.Range("A2:B9").Select()
